$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 222 (shifts the existing 222..325 down to 223..326),
# then populate it with the new weekly record (2022-10-12).
$ws.Rows.Item(222).Insert()
$ws.Range("A222").Value = 5
$ws.Range("B222").Value = "Macroferia Regional de Talca"
$ws.Range("C222").Value = "Maule"
$ws.Range("D222").Value = 44846
$ws.Range("E222").Value = 7
$ws.Range("F222").Value = 100112009
$ws.Range("G222").Value = "Acelga"
$ws.Range("H222").Value = "Sin especificar"
$ws.Range("I222").Value = "Primera"
$ws.Range("J222").Value = 500
$ws.Range("K222").Value = 2500
$ws.Range("L222").Value = 2500
$ws.Range("M222").Value = 2500
$ws.Range("N222").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O222").Value = "Región del Maule"
$ws.Range("P222").Value = 625
$ws.Range("Q222").Value = 4
$ws.Range("R222").Value = "Hortaliza"

# Insert a second new row at 228 (shifts 228..326 down to 229..327),
# then populate it with the other new weekly record (2022-10-11).
$ws.Rows.Item(228).Insert()
$ws.Range("A228").Value = 5
$ws.Range("B228").Value = "Macroferia Regional de Talca"
$ws.Range("C228").Value = "Maule"
$ws.Range("D228").Value = 44845
$ws.Range("E228").Value = 7
$ws.Range("F228").Value = 100112009
$ws.Range("G228").Value = "Acelga"
$ws.Range("H228").Value = "Sin especificar"
$ws.Range("I228").Value = "Primera"
$ws.Range("J228").Value = 500
$ws.Range("K228").Value = 2500
$ws.Range("L228").Value = 2500
$ws.Range("M228").Value = 2500
$ws.Range("N228").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O228").Value = "Región del Maule"
$ws.Range("P228").Value = 625
$ws.Range("Q228").Value = 4
$ws.Range("R228").Value = "Hortaliza"
